# Update "想去人数" (column F) counts on the 展览, 演出 and 全部类型 sheets.
# This mirrors the data refresh performed by the gh-pages publishing job
# (commit "Update gh-pages to output generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 35
$ws1.Range("F3").Value = 71
$ws1.Range("F4").Value = 83
$ws1.Range("F5").Value = 961
$ws1.Range("F6").Value = 357
$ws1.Range("F8").Value = 554
$ws1.Range("F9").Value = 1437
$ws1.Range("F11").Value = 1328
$ws1.Range("F12").Value = 2991
$ws1.Range("F13").Value = 390
$ws1.Range("F14").Value = 1603
$ws1.Range("F16").Value = 785
$ws1.Range("F17").Value = 234
$ws1.Range("F18").Value = 1364
$ws1.Range("F19").Value = 262
$ws1.Range("F21").Value = 1115
$ws1.Range("F22").Value = 395
$ws1.Range("F23").Value = 3457
$ws1.Range("F24").Value = 675
$ws1.Range("F26").Value = 1529

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 168
$ws2.Range("F5").Value = 24
$ws2.Range("F7").Value = 48
$ws2.Range("F12").Value = 77

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 35
$ws4.Range("F5").Value = 71
$ws4.Range("F8").Value = 168
$ws4.Range("F9").Value = 24
$ws4.Range("F11").Value = 48
$ws4.Range("F12").Value = 83
$ws4.Range("F15").Value = 961
$ws4.Range("F16").Value = 357
$ws4.Range("F18").Value = 554
$ws4.Range("F19").Value = 1437
$ws4.Range("F21").Value = 1328
$ws4.Range("F22").Value = 2991
$ws4.Range("F23").Value = 390
$ws4.Range("F24").Value = 1603
$ws4.Range("F26").Value = 785
$ws4.Range("F27").Value = 234
$ws4.Range("F28").Value = 1364
$ws4.Range("F29").Value = 262
$ws4.Range("F33").Value = 1115
$ws4.Range("F34").Value = 395
$ws4.Range("F35").Value = 3457
$ws4.Range("F36").Value = 675
$ws4.Range("F38").Value = 1529
$ws4.Range("F39").Value = 77
